# Apply odds/statistics updates to the FlashScore weekly games sheet (2024-10-28).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.25
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 1.8
$ws.Range("L2").Value = 13
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("U2").Value = 3.4
$ws.Range("V2").Value = 1.3
$ws.Range("X2").Value = 4.5
$ws.Range("Z2").Value = 6.5
$ws.Range("AA2").Value = 17
$ws.Range("AD2").Value = 11
$ws.Range("AI2").Value = 67
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 17
$ws.Range("AU2").Value = 15
$ws.Range("AW2").Value = 13
$ws.Range("AY2").Value = 81

# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 3.5
$ws.Range("AA5").Value = 19
$ws.Range("AI5").Value = 15
$ws.Range("AK5").Value = 29
$ws.Range("AN5").Value = 4.5
$ws.Range("AO5").Value = 13
$ws.Range("AS5").Value = 151

# Row 6
$ws.Range("G6").Value = 10.75
$ws.Range("H6").Value = 5.9
$ws.Range("I6").Value = 1.18
$ws.Range("J6").Value = 8.5
$ws.Range("K6").Value = 2.9
$ws.Range("L6").Value = 1.5
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 19.5
$ws.Range("P6").Value = 5.6
$ws.Range("Q6").Value = 1.33
$ws.Range("R6").Value = 2.73
$ws.Range("U6").Value = 1.92
$ws.Range("V6").Value = 1.84
$ws.Range("W6").Value = 29
$ws.Range("X6").Value = 75
$ws.Range("Y6").Value = 29
$ws.Range("AA6").Value = 100
$ws.Range("AB6").Value = 75
$ws.Range("AC6").Value = 18.5
$ws.Range("AD6").Value = 11.25
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 8
$ws.Range("AI6").Value = 6
$ws.Range("AK6").Value = 6.3
$ws.Range("AL6").Value = 8.75
$ws.Range("AM6").Value = 21
$ws.Range("AN6").Value = 12
$ws.Range("AO6").Value = 65
$ws.Range("AP6").Value = 50
$ws.Range("AQ6").Value = 500
$ws.Range("AT6").Value = 4.1
$ws.Range("AX6").Value = 4.75
$ws.Range("AZ6").Value = 10
$ws.Range("BA6").Value = 28

# Row 9
$ws.Range("G9").Value = 2.25
$ws.Range("I9").Value = 2.8
$ws.Range("Z9").Value = 21
$ws.Range("AC9").Value = 10
$ws.Range("AE9").Value = 15
$ws.Range("AX9").Value = 17

# Row 11
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.5

# Row 12
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93

# Row 13
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 2.75
$ws.Range("I13").Value = 3.1
$ws.Range("J13").Value = 3.5
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 1.14
$ws.Range("N13").Value = 5.5
$ws.Range("O13").Value = 1.62
$ws.Range("P13").Value = 2.2
$ws.Range("Q13").Value = 3.1
$ws.Range("R13").Value = 1.36
$ws.Range("S13").Value = 1.67
$ws.Range("T13").Value = 2.1
$ws.Range("W13").Value = 6
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 12
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 29
$ws.Range("AD13").Value = 5.5
$ws.Range("AI13").Value = 13
$ws.Range("AK13").Value = 34
$ws.Range("AO13").Value = 17
$ws.Range("AT13").Value = 2.1
$ws.Range("AW13").Value = 4.75
$ws.Range("AY13").Value = 34
$ws.Range("AZ13").Value = 67

# Row 14
$ws.Range("G14").Value = 1.6
$ws.Range("I14").Value = 5.25
$ws.Range("Q14").Value = 1.7
$ws.Range("R14").Value = 2.1
$ws.Range("AA14").Value = 12
$ws.Range("AQ14").Value = 23

# Row 17
$ws.Range("G17").Value = 1.8
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 3.9
$ws.Range("J17").Value = 2.3
$ws.Range("K17").Value = 2.63
$ws.Range("L17").Value = 3.75
$ws.Range("X17").Value = 13
$ws.Range("AI17").Value = 26
$ws.Range("AM17").Value = 23
$ws.Range("AN17").Value = 4.5
$ws.Range("AP17").Value = 13
$ws.Range("AQ17").Value = 26
$ws.Range("AY17").Value = 19
$ws.Range("BB17").Value = 81
$ws.Range("BC17").Value = 201

# Row 18
$ws.Range("G18").Value = 1.8
$ws.Range("H18").Value = 3.5
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 2.4
$ws.Range("K18").Value = 2.2
$ws.Range("L18").Value = 4.5
$ws.Range("M18").Value = 1.05
$ws.Range("N18").Value = 11
$ws.Range("O18").Value = 1.29
$ws.Range("P18").Value = 3.5
$ws.Range("Q18").Value = 1.9
$ws.Range("R18").Value = 1.9
$ws.Range("S18").Value = 1.36
$ws.Range("T18").Value = 3
$ws.Range("U18").Value = 1.8
$ws.Range("V18").Value = 1.91
$ws.Range("W18").Value = 7.5
$ws.Range("X18").Value = 8.5
$ws.Range("AA18").Value = 15
$ws.Range("AB18").Value = 26
$ws.Range("AC18").Value = 11
$ws.Range("AD18").Value = 7
$ws.Range("AF18").Value = 51
$ws.Range("AG18").Value = 251
$ws.Range("AH18").Value = 12
$ws.Range("AI18").Value = 21
$ws.Range("AJ18").Value = 15
$ws.Range("AL18").Value = 34
$ws.Range("AM18").Value = 41
$ws.Range("AN18").Value = 3.75
$ws.Range("AP18").Value = 21
$ws.Range("AQ18").Value = 34
$ws.Range("AR18").Value = 51
$ws.Range("AS18").Value = 126
$ws.Range("AT18").Value = 3
$ws.Range("AU18").Value = 8
$ws.Range("AX18").Value = 23
$ws.Range("AY18").Value = 29
$ws.Range("AZ18").Value = 81
$ws.Range("BA18").Value = 101
$ws.Range("BB18").Value = 201
$ws.Range("BD18").Value = 126
